$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.004.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.063.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.056.32"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.20"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.558.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.005.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.066.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.691"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.04"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.81%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "58.90"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.41"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "476.46"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.250.57"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0393"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0787"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.07"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.11"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0516"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.01%  "
